$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# Update the "so cong thuc hien" (work count) values for rows 11-19, column E
$ws.Range("E11").Value = 5
$ws.Range("E12").Value = 12
$ws.Range("E13").Value = 24
$ws.Range("E14").Value = 12
$ws.Range("E15").Value = 24
$ws.Range("E16").Value = 20
$ws.Range("E17").Value = 80
$ws.Range("E18").Value = 24
$ws.Range("E19").Value = 24

# Update the view state: scrolled so row 7 is at the top, and E17 is the active cell
$ws.Range("E17").Select()
$excel.ActiveWindow.ScrollRow = 7
